# Applies the "Updated cryptos list" GitHub Actions refresh commit.
# Source data is scraped crypto market info; Price (D) and Volume(1h) (E)
# are plain text cells (values like "70.683.62" aren't valid Excel numbers,
# so the sheet stores them as text). A couple of rows (OKB / FirstDigitalUSD,
# and TheGraph / ThetaToken) also swapped positions in this refresh, which
# shows up as their Coin/Link/Price/Volume cells changing together.
#
# Some refreshed price values (e.g. "618.99", "9.20", "1.00") *do* parse as
# plain numbers. Typing such a value into a general-formatted cell would
# normally have Excel coerce it to a Number (losing the original text
# formatting/trailing zeros), so for those cells we force the cell to Text
# format first, matching the source workbook's inline-string cell type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.697.26'
$ws.Range("E2").Value = '  -0.61%  '
$ws.Range("D3").Value = '3.520.32'
$ws.Range("E3").Value = '  -2.41%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '618.99'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.03'
$ws.Range("E6").Value = '  -1.02%  '
$ws.Range("E7").Value = '  -1.16%  '
$ws.Range("D8").Value = '3.516.00'
$ws.Range("E8").Value = '  -2.30%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  -0.12%  '
$ws.Range("E10").Value = '  -2.28%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.09'
$ws.Range("E11").Value = '  -4.77%  '
$ws.Range("E12").Value = '  -1.04%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '46.42'
$ws.Range("E13").Value = '  -1.92%  '
$ws.Range("E14").Value = '  -1.91%  '
$ws.Range("D15").Value = '4.091.70'
$ws.Range("E15").Value = '  -2.26%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.37'
$ws.Range("E16").Value = '  -1.45%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '608.74'
$ws.Range("E17").Value = '  -1.97%  '
$ws.Range("D18").Value = '3.527.05'
$ws.Range("E18").Value = '  -2.14%  '
$ws.Range("D19").Value = '70.784.95'
$ws.Range("E19").Value = '  -0.56%  '
$ws.Range("E20").Value = '  +1.15%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.73'
$ws.Range("E21").Value = '  +0.89%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.884'
$ws.Range("E22").Value = '  -1.20%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.20'
$ws.Range("E23").Value = '  -1.45%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '15.62'
$ws.Range("E24").Value = '  -3.55%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '97.83'
$ws.Range("E25").Value = '  -0.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.72'
$ws.Range("E26").Value = '  -2.57%  '
$ws.Range("E27").Value = '  -0.01%  '
$ws.Range("E28").Value = '  -4.21%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.66'
$ws.Range("E29").Value = '  -1.96%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.07'
$ws.Range("E30").Value = '  -3.17%  '
$ws.Range("E31").Value = '  -3.67%  '
$ws.Range("E32").Value = '  -5.65%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.30'
$ws.Range("E33").Value = '  -1.27%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '641.13'
$ws.Range("E34").Value = '  +1.96%  '
$ws.Range("E35").Value = '  -6.68%  '
$ws.Range("E37").Value = '  -1.21%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0486'
$ws.Range("E38").Value = '  -0.20%  '
$ws.Range("E39").Value = '  -9.03%  '
$ws.Range("B40").Value = 'OKB'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '56.72'
$ws.Range("E40").Value = '  -1.71%  '
$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.37%  '
$ws.Range("E42").Value = '  -0.75%  '
$ws.Range("D43").Value = '3.351.75'
$ws.Range("E43").Value = '  -1.69%  '
$ws.Range("D44").Value = '0.0₃0719'
$ws.Range("E44").Value = '  -0.68%  '
$ws.Range("B45").Value = 'ThetaToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.96'
$ws.Range("E45").Value = '  -1.72%  '
$ws.Range("B46").Value = 'TheGraph'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.312'
$ws.Range("E46").Value = '  -5.05%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '31.80'
$ws.Range("E47").Value = '  -4.24%  '
$ws.Range("E48").Value = '  -7.11%  '
$ws.Range("E49").Value = '  -0.68%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '134.52'
$ws.Range("E50").Value = '  +1.23%  '
$ws.Range("E51").Value = '  -0.02%  '
